$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while forcing it to remain plain text
# (prevents Excel from auto-converting numeric-looking strings like "6.719"
# or "0.000009072" into Number/Date types), then restores the original style
# so no stray formatting/quote-prefix style is left behind.
function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = $origStyle
}

# Row 2
Set-TextValue "D2" '28.247.61'
Set-TextValue "E2" '  -2.36%  '

# Row 3
Set-TextValue "D3" '1.866.70'
Set-TextValue "E3" '  -1.81%  '

# Row 4
Set-TextValue "D4" '1.005'
Set-TextValue "E4" '  +0.14%  '

# Row 5
Set-TextValue "D5" '318.92'
Set-TextValue "E5" '  -1.75%  '

# Row 6
Set-TextValue "D6" '1.003'
Set-TextValue "E6" '  +0.13%  '

# Row 7
Set-TextValue "D7" '0.4378'
Set-TextValue "E7" '  -4.54%  '

# Row 8
Set-TextValue "D8" '0.3700'
Set-TextValue "E8" '  -3.09%  '

# Row 9
Set-TextValue "D9" '0.07511'
Set-TextValue "E9" '  -2.55%  '

# Row 10
Set-TextValue "D10" '0.9387'
Set-TextValue "E10" '  -3.84%  '

# Row 11
Set-TextValue "D11" '21.38'
Set-TextValue "E11" '  -3.03%  '

# Row 12
Set-TextValue "D12" '1.893.55'
Set-TextValue "E12" '  -0.17%  '

# Row 13
Set-TextValue "D13" '6.719'
Set-TextValue "E13" '  -3.12%  '

# Row 14
Set-TextValue "D14" '5.442'
Set-TextValue "E14" '  -3.60%  '

# Row 15
Set-TextValue "D15" '0.06866'
Set-TextValue "E15" '  -2.62%  '

# Row 16
Set-TextValue "D16" '1.005'
Set-TextValue "E16" '  +0.24%  '

# Row 17
Set-TextValue "D17" '82.20'
Set-TextValue "E17" '  -1.59%  '

# Row 18
Set-TextValue "D18" '0.000009072'

# Row 20
Set-TextValue "D20" '15.93'
Set-TextValue "E20" '  -4.22%  '

# Row 21
Set-TextValue "D21" '28.256.42'
Set-TextValue "E21" '  -2.24%  '

# Row 22
Set-TextValue "D22" '5.131'
Set-TextValue "E22" '  -3.12%  '

# Row 23
Set-TextValue "D23" '10.77'
Set-TextValue "E23" '  -0.72%  '

# Row 24
Set-TextValue "D24" '2.155.04'
Set-TextValue "E24" '  +1.01%  '

# Row 25
Set-TextValue "D25" '2.028'
Set-TextValue "E25" '  -3.10%  '

# Row 26
Set-TextValue "D26" '154.94'
Set-TextValue "E26" '  -1.89%  '

# Row 27
Set-TextValue "D27" '18.42'
Set-TextValue "E27" '  -3.04%  '

# Row 28
Set-TextValue "D28" '5.309'
Set-TextValue "E28" '  -5.74%  '

# Row 29
Set-TextValue "D29" '113.94'
Set-TextValue "E29" '  -2.92%  '

# Row 30
Set-TextValue "D30" '1.727'
Set-TextValue "E30" '  -5.94%  '

# Row 31
Set-TextValue "D31" '0.09050'
Set-TextValue "E31" '  -2.17%  '

# Row 32
Set-TextValue "D32" '0.7958'
Set-TextValue "E32" '  -7.78%  '

# Row 33
Set-TextValue "D33" '4.833'
Set-TextValue "E33" '  -4.89%  '

# Row 34
Set-TextValue "D34" '1.168'
Set-TextValue "E34" '  -5.84%  '

# Row 35
Set-TextValue "D35" '2.961'
Set-TextValue "E35" '  -1.19%  '

# Row 36
Set-TextValue "D36" '1.003'
Set-TextValue "E36" '  +0.12%  '

# Row 37
Set-TextValue "D37" '1.118'
Set-TextValue "E37" '  -2.46%  '

# Row 38
Set-TextValue "D38" '0.05426'
Set-TextValue "E38" '  -5.02%  '

# Row 39
Set-TextValue "D39" '0.01954'
Set-TextValue "E39" '  -4.24%  '

# Row 40
Set-TextValue "D40" '2.943'
Set-TextValue "E40" '  +6.40%  '

# Row 41
Set-TextValue "D41" '7.125'
Set-TextValue "E41" '  -3.56%  '

# Row 42
Set-TextValue "D42" '0.5247'
Set-TextValue "E42" '  -4.33%  '

# Row 43
Set-TextValue "D43" '0.1671'
Set-TextValue "E43" '  -4.53%  '

# Row 44
Set-TextValue "D44" '8.690'
Set-TextValue "E44" '  -6.25%  '

# Row 45
Set-TextValue "D45" '0.06762'
Set-TextValue "E45" '  -0.67%  '

# Row 46
$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-TextValue "D46" '0.4874'
Set-TextValue "E46" '  -5.50%  '

# Row 47
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue "D47" '10.58'
Set-TextValue "E47" '  -5.72%  '

# Row 48
$ws.Range("B48").Value = 'RenderToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue "D48" '1.986'
Set-TextValue "E48" '  -3.71%  '

# Row 49
Set-TextValue "D49" '107.93'
Set-TextValue "E49" '  -2.03%  '

# Row 50
$ws.Range("B50").Value = 'PaxDollar'
$ws.Range("C50").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextValue "D50" '1.003'
Set-TextValue "E50" '  +0.11%  '

# Row 51
Set-TextValue "D51" '1.675'
Set-TextValue "E51" '  -5.48%  '

